# Sigi-9: Find rows by headers
# Adds a second data column (C) to the SIGIDOC CELLS ENG test fixture so the
# "find rows by header" lookup logic can be exercised against more than one
# data column. Mirrors column B's values for a couple of rows and adds a few
# new sample values used by the new tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C needs roughly the same kind of width as column B so the longer
# sample strings are readable.
$ws.Columns("C").ColumnWidth = 19.3

# SEAL ID (row 1) / TYPE (row 2) - simple duplicate / companion example values
$ws.Range("C1").Value = 98
$ws.Range("C2").Value = "Seal"

# Title / Editor forename / Editor surname (rows 57-59)
$ws.Range("C57").Value = "This is a seal title"
$ws.Range("C58").Value = "Steve"
$ws.Range("C59").Value = "Shaw"

# EDITION(S) (row 45) - second edition reference
$ws.Range("C45").Value = "Ivanov, 2017, 32, p.47"

# Filename (row 60)
$ws.Range("C60").Value = "TM_98"

# Sequence (row 61) - stored as text (zero padded) rather than a number, so
# force a text number format before assigning the values.
$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = "0099"

$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "0098"

# Leave the selection near the new data, like the source workbook.
$ws.Range("C62").Select()
